# Actualización automática hashcode sáb may 25 02:08:53 CEST 2019
# Updates hash values in column B of the hashcode.csv sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B15").Value = "9121218654e7cdae952b187a3f158493"
$ws.Range("B16").Value = "a5c13d80c47c332fa1ee0c687c1511d8"
$ws.Range("B35").Value = "1767a306e50a67b63060a319c55ca1b2"
$ws.Range("B49").Value = "9702ab74fbe3ecfe7d35fad3d39d57c4"
$ws.Range("B154").Value = "f667232bc9fcec8ed198a5c752d39832"
$ws.Range("B159").Value = "57bab5ae699de9a5a6f6cef6d0a1d855"
$ws.Range("B198").Value = "c8eb21a68ef50809b9b43c6152f25989"
$ws.Range("B227").Value = "8a0aa50bfbd30ae778a964bc3ccef7d8"
$ws.Range("B232").Value = "d53ce3f95a81937bf61b93ea482d1bb9"
$ws.Range("B280").Value = "492e2c8558dd97a9d0b23fd3851b3587"
$ws.Range("B281").Value = "01bbae34a33430e69b19f9960d490a10"
$ws.Range("B299").Value = "ca06a29ddf84c1012ce23445464311d1"
$ws.Range("B339").Value = "df0e2d182d7efbef33009bc513503800"
$ws.Range("B350").Value = "205045de71ccf4d8ebb7043be63d7d1e"
$ws.Range("B358").Value = "fa7d097d702a3fb4c02ab0e6ec24568c"
$ws.Range("B397").Value = "e899843e8de1d189c9c71a6969ab9f97"
$ws.Range("B419").Value = "930e9bd628ccd09c643cd2b4a4b8cfad"
$ws.Range("B424").Value = "c3d15ba386f49a4a89cff768392ffa95"
$ws.Range("B451").Value = "e978f9e6327940d2114367ce417a223d"
$ws.Range("B460").Value = "feb5208e183874b5f77091b49b34c766"
$ws.Range("B465").Value = "fd97ec2bfb2c11dc87f2ba81f8bca5a3"
$ws.Range("B478").Value = "d8912b074d0b14438de67954956522e7"
$ws.Range("B500").Value = "d23276c9d7611ab2a179a914c1ffc24b"
$ws.Range("B502").Value = "f4cb37c7b35fcd483f14ef6de32a5d79"
$ws.Range("B517").Value = "139ed9d1975ba0e500b517e51e222364"
$ws.Range("B569").Value = "0e540c3ca72bac2a30b6f2007ef86776"
$ws.Range("B616").Value = "be189cad32e337fbd4375c17fcde0b8a"
$ws.Range("B627").Value = "5c79a81ffd9cdd47fd6b6a98597dacc4"
$ws.Range("B665").Value = "beca78e553b4e90b954c318ee3b77654"
$ws.Range("B712").Value = "f5c07954d5e36d9a67fc8c20c5548bcb"
$ws.Range("B726").Value = "46b26a82f804b99a37557d8cd64c692c"
$ws.Range("B731").Value = "a561d1bf4aefcf39e61e1863b8147b44"
$ws.Range("B733").Value = "bfc7424cfab780ad37061f71bb9397ad"
$ws.Range("B740").Value = "d4374f0fa39c6f7edfbd28cca214f2b8"
$ws.Range("B756").Value = "a833e19224e0d52ed7ff59b2093d743d"
$ws.Range("B802").Value = "11e6135d92906710ca6283d07f1d1add"
$ws.Range("B811").Value = "7ab1965ef13b28203fff68fba8a07b70"
$ws.Range("B819").Value = "26eeb802ed9231e661d19f47fc07dc10"
$ws.Range("B839").Value = "838e687b650fda7a6da60c9e4c56a4be"
$ws.Range("B846").Value = "da70563953f6e5c1d4a1aab0bbe1d7e0"
$ws.Range("B848").Value = "661c7a2286dd8390bd5f9d2ff11d671b"
$ws.Range("B888").Value = "d0b7821b6e6a30385eb91b05fb4adc9f"
$ws.Range("B911").Value = "cba30d7950a13a0c0967661dd8f1ded3"
$ws.Range("B937").Value = "7bc9e2f9a2e884246334334a01f323a4"
$ws.Range("B951").Value = "3f574683856d8cc29639b08f7ab41e07"